$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Jam2"
$ws.Range("C2").Value = "F11r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 33.24999533333333
$ws.Range("H2").Value = 99.74998599999999
$ws.Range("I2").Value = 0.5673360890306117
$ws.Range("J2").Value = 0.5673360890306117
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 20.10858766666667
$ws.Range("N2").Value = 60.325763
$ws.Range("O2").Value = 0.9234633461941997
$ws.Range("P2").Value = 0.9234633461941997
$ws.Range("Q2").Value = 668.6104460765908
$ws.Range("R2").Value = 6017.494014689318
$ws.Range("S2").Value = 0.5239140831929391
$ws.Range("T2").Value = 0.5239140831929391

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Jam2"
$ws.Range("C3").Value = "F11r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 33.24999533333333
$ws.Range("H3").Value = 99.74998599999999
$ws.Range("I3").Value = 0.5673360890306117
$ws.Range("J3").Value = 0.5673360890306117
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.2128156666666667
$ws.Range("N3").Value = 0.638447
$ws.Range("O3").Value = 0.009773310334883756
$ws.Range("P3").Value = 0.009773310334883756
$ws.Range("Q3").Value = 7.076119923526888
$ws.Range("R3").Value = 63.68507931174199
$ws.Range("S3").Value = 0.005544751662275409
$ws.Range("T3").Value = 0.005544751662275409

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Jam2"
$ws.Range("C4").Value = "F11r"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 33.24999533333333
$ws.Range("H4").Value = 99.74998599999999
$ws.Range("I4").Value = 0.5673360890306117
$ws.Range("J4").Value = 0.5673360890306117
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.453784333333333
$ws.Range("N4").Value = 4.361352999999999
$ws.Range("O4").Value = 0.06676334347091657
$ws.Range("P4").Value = 0.06676334347091657
$ws.Range("Q4").Value = 48.33832229900644
$ws.Range("R4").Value = 435.0449006910579
$ws.Range("S4").Value = 0.03787725417539724
$ws.Range("T4").Value = 0.03787725417539724

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Jam2"
$ws.Range("C5").Value = "F11r"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 23.30243966666667
$ws.Range("H5").Value = 69.907319
$ws.Range("I5").Value = 0.3976035140102714
$ws.Range("J5").Value = 0.3976035140102714
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 20.10858766666667
$ws.Range("N5").Value = 60.325763
$ws.Range("O5").Value = 0.9234633461941997
$ws.Range("P5").Value = 0.9234633461941997
$ws.Range("Q5").Value = 468.5791508843775
$ws.Range("R5").Value = 4217.212357959397
$ws.Range("S5").Value = 0.3671722715064976
$ws.Range("T5").Value = 0.3671722715064976

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Jam2"
$ws.Range("C6").Value = "F11r"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 23.30243966666667
$ws.Range("H6").Value = 69.907319
$ws.Range("I6").Value = 0.3976035140102714
$ws.Range("J6").Value = 0.3976035140102714
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.2128156666666667
$ws.Range("N6").Value = 0.638447
$ws.Range("O6").Value = 0.009773310334883756
$ws.Range("P6").Value = 0.009773310334883756
$ws.Range("Q6").Value = 4.959124232621445
$ws.Range("R6").Value = 44.632118093593
$ws.Range("S6").Value = 0.003885902532662684
$ws.Range("T6").Value = 0.003885902532662684

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Jam2"
$ws.Range("C7").Value = "F11r"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 23.30243966666667
$ws.Range("H7").Value = 69.907319
$ws.Range("I7").Value = 0.3976035140102714
$ws.Range("J7").Value = 0.3976035140102714
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.453784333333333
$ws.Range("N7").Value = 4.361352999999999
$ws.Range("O7").Value = 0.06676334347091657
$ws.Range("P7").Value = 0.06676334347091657
$ws.Range("Q7").Value = 33.87672171584522
$ws.Range("R7").Value = 304.8904954426069
$ws.Range("S7").Value = 0.02654533997111114
$ws.Range("T7").Value = 0.02654533997111114

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Jam2"
$ws.Range("C8").Value = "F11r"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.054792666666667
$ws.Range("H8").Value = 6.164378
$ws.Range("I8").Value = 0.03506039695911681
$ws.Range("J8").Value = 0.03506039695911681
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 20.10858766666667
$ws.Range("N8").Value = 60.325763
$ws.Range("O8").Value = 0.9234633461941997
$ws.Range("P8").Value = 0.9234633461941997
$ws.Range("Q8").Value = 41.31897847449045
$ws.Range("R8").Value = 371.870806270414
$ws.Range("S8").Value = 0.03237699149476295
$ws.Range("T8").Value = 0.03237699149476295

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Jam2"
$ws.Range("C9").Value = "F11r"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.054792666666667
$ws.Range("H9").Value = 6.164378
$ws.Range("I9").Value = 0.03506039695911681
$ws.Range("J9").Value = 0.03506039695911681
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2128156666666667
$ws.Range("N9").Value = 0.638447
$ws.Range("O9").Value = 0.009773310334883756
$ws.Range("P9").Value = 0.009773310334883756
$ws.Range("Q9").Value = 0.4372920712184444
$ws.Range("R9").Value = 3.935628640966
$ws.Range("S9").Value = 0.0003426561399456634
$ws.Range("T9").Value = 0.0003426561399456633

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Jam2"
$ws.Range("C10").Value = "F11r"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.054792666666667
$ws.Range("H10").Value = 6.164378
$ws.Range("I10").Value = 0.03506039695911681
$ws.Range("J10").Value = 0.03506039695911681
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.453784333333333
$ws.Range("N10").Value = 4.361352999999999
$ws.Range("O10").Value = 0.06676334347091657
$ws.Range("P10").Value = 0.06676334347091657
$ws.Range("Q10").Value = 2.987225387048222
$ws.Range("R10").Value = 26.885028483434
$ws.Range("S10").Value = 0.002340749324408195
$ws.Range("T10").Value = 0.002340749324408194
